$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 10
$ws.Range("A3").Value = 20
$ws.Range("A4").Value = 30
$ws.Range("A5").Value = 40

$ws.Range("D5").Select()
